# Insert two new columns (Q and R) for "default_count" and "default_value",
# shifting the existing most_frequent_value / memory_consumed_bytes /
# pattern_count / patterns columns two places to the right.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("Q:R").Insert()

# New header cells.
$ws.Range("Q1").Value = "default_count"
$ws.Range("R1").Value = "default_value"

# New data columns for every data row (2-21).
$ws.Range("Q2:Q21").Value = 0
$ws.Range("R2:R21").Value = "<Unspecified>"

# The recomputed "most_frequent_value" column (now S) changed for several
# rows as part of this profiling run.
$ws.Range("S2").Value = ""
$ws.Range("S3").Value = "Published"
$ws.Range("S4").Value = "Murarrie station"
$ws.Range("S5").Value = "Grey Street, South Brisbane, Brisbane Central"
$ws.Range("S6").Value = "-27.4818"
$ws.Range("S7").Value = "153.029"
$ws.Range("S8").Value = ""
$ws.Range("S9").Value = "1"
$ws.Range("S10").Value = "Standard retailer"
$ws.Range("S11").Value = ""
$ws.Range("S12").Value = "go card"
$ws.Range("S13").Value = ""
$ws.Range("S14").Value = ""
$ws.Range("S15").Value = "Buy"
$ws.Range("S16").Value = ""
$ws.Range("S17").Value = ""
$ws.Range("S18").Value = ""
$ws.Range("S19").Value = ""
$ws.Range("S20").Value = ""
$ws.Range("S21").Value = ""
